$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.674.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.885.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5128"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3959"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08447"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.118"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.84"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.313"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.886.49"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.57"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.310"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001111"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06744"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.81"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.994"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.698.75"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.257"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.101.57"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.81"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.397"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.19"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1056"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.056"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.833"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.623"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02476"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06557"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2201"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.986"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.267"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.126"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6499"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.23"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.008"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6100"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.712"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.053"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.226"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.94"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.200"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.43%  "
